# Corrected name of the test class for task 2 criteria
# (CheckDataTimeStringTest.java -> DateTimeCheckerTest.java)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$oldSegment = "in CheckDataTimeStringTest.java (2"
$newSegment = "in  DateTimeCheckerTest.java  (2"

$cell = $ws.Range("D10")
$newVal = $cell.Value2.Replace($oldSegment, $newSegment)
$cell.Value = $newVal

# Move the active selection to D10, matching the saved workbook state
$ws.Range("D10").Select()
